$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column E (Level) values for rows 2-18
$ws.Range("E2:E9").Value = 0
$ws.Range("E10:E14").Value = 1
$ws.Range("E15:E18").Value = 2

# Update the selection to U5:U18 with active cell U5
$ws.Range("U5:U18").Select()
